# In the "Recorded By" column (G) of the "Session Analysis Results" sheet,
# every cell that reads "dnasr281@gmail.com, System" has the two names
# reordered to "System, dnasr281@gmail.com".
#
# Using Range.Replace (rather than iterating/reading every cell) so that
# blank cells in the column are left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$col = $ws.Columns.Item(7)  # Column G = "Recorded By"

# xlWhole=1, xlByRows=1 (unused by whole-cell match), not match case, no match byte
[void]$col.Replace($oldValue, $newValue, 1, 1, $false, $false, $false)
